{"js": "// Add a \"Meta description\" paragraph right after the title, and turn the\n// old \"Play Book of Helios Free Slot Game Online\" / \"Read our review...\"\n// pair at the bottom of the doc into a single DALLE image-prompt paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// --- 1. Insert a new \"Meta description: ...\" paragraph right after the\n// Heading1 title paragraph, before \"Book of Helios: Shine Bright...\" ---\nconst titlePara = paragraphs.items[0];\nconst metaLabel = \"Meta description\";\nconst metaRest =\n  \": Read our review of Book of Helios, a Betsoft slot game themed around \" +\n  \"Helios. Play for free and access the free spins round with buy feature option.\";\n\nconst metaPara = titlePara.insertParagraph(metaLabel + metaRest, Word.InsertLocation.after);\n// The body paragraphs in this document use the \"Normal\" style (the new\n// paragraph would otherwise inherit the Heading1 style from its anchor).\nmetaPara.style = \"Normal\";\nawait context.sync();\n\n// Bold just the \"Meta description\" label, leaving the rest of the sentence\n// in regular formatting.\nconst metaLabelRange = metaPara.search(metaLabel, { matchCase: true });\nmetaLabelRange.load(\"text\");\nawait context.sync();\nmetaLabelRange.items[0].font.bold = true;\nawait context.sync();\n\n// --- 2 & 3. At the bottom of the document: drop the duplicated bold\n// \"Play Book of Helios Free Slot Game Online\" paragraph, and replace the\n// text of the final (italic) paragraph with the DALLE prompt. ---\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst count = paragraphs.items.length;\nconst duplicateTitlePara = paragraphs.items[count - 2];\nconst reviewPara = paragraphs.items[count - 1];\n\nduplicateTitlePara.delete();\nawait context.sync();\n\nconst oldReviewText =\n  \"Read our review of Book of Helios, a Betsoft slot game themed around Helios. \" +\n  \"Play for free and access the free spins round with buy feature option.\";\nconst dallePrompt =\n  'DALLE, please create a cartoon-style feature image for the game \"Book of Helios\" ' +\n  \"that features a happy Maya warrior with glasses. The image should be captivating \" +\n  \"and exciting, with the Maya warrior shown holding a copy of the book with Helios \" +\n  \"on the cover. The warrior should have a big smile and be surrounded by golden rays \" +\n  \"of sunlight. The image should also include other symbols from the game, such as the \" +\n  \"Book of Helios symbol and the expandable symbol. Make sure the overall color scheme \" +\n  \"is bright and eye-catching, and that the image is of high-quality. Thanks!\";\n\nconst reviewTextRange = reviewPara.search(oldReviewText, { matchCase: true });\nreviewTextRange.load(\"text\");\nawait context.sync();\n// Replacing just the matched range (instead of the whole paragraph) keeps\n// the paragraph's existing run/formatting structure (italic) intact.\nreviewTextRange.items[0].insertText(dallePrompt, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Add a \"Meta description\" paragraph right after the title, and turn the\n# old \"Play Book of Helios Free Slot Game Online\" / \"Read our review...\"\n# pair at the bottom of the doc into a single DALLE image-prompt paragraph.\n\n$d = $word.ActiveDocument\n\n# --- 1. Insert a new \"Meta description: ...\" paragraph right after the\n# Heading1 title paragraph, before \"Book of Helios: Shine Bright...\" ---\n$metaLabel = \"Meta description\"\n$metaRest = \": Read our review of Book of Helios, a Betsoft slot game themed around Helios. Play for free and access the free spins round with buy feature option.\"\n\n$titlePara = $d.Paragraphs(1)\n$titlePara.Range.InsertParagraphAfter() | Out-Null\n\n$metaPara = $d.Paragraphs(2)\n# The body paragraphs in this document use the \"Normal\" style (the new\n# paragraph would otherwise inherit the Heading1 style from its anchor).\n$metaPara.Style = \"Normal\"\n$metaPara.Range.Text = $metaLabel + $metaRest\n\n# Bold just the \"Meta description\" label, leaving the rest of the sentence\n# in regular formatting.\n$metaLabelRange = $metaPara.Range.Duplicate\n$metaLabelRange.Find.Execute($metaLabel) | Out-Null\n$metaLabelRange.Bold = 1\n\n# --- 2 & 3. At the bottom of the document: drop the duplicated bold\n# \"Play Book of Helios Free Slot Game Online\" paragraph, and replace the\n# text of the final (italic) paragraph with the DALLE prompt. ---\n$count = $d.Paragraphs.Count\n$duplicateTitlePara = $d.Paragraphs($count - 1)\n$duplicateTitlePara.Range.Delete() | Out-Null\n\n$newCount = $d.Paragraphs.Count\n$reviewPara = $d.Paragraphs($newCount)\n\n$oldReviewText = \"Read our review of Book of Helios, a Betsoft slot game themed around Helios. Play for free and access the free spins round with buy feature option.\"\n$dallePrompt = 'DALLE, please create a cartoon-style feature image for the game \"Book of Helios\" that features a happy Maya warrior with glasses. The image should be captivating and exciting, with the Maya warrior shown holding a copy of the book with Helios on the cover. The warrior should have a big smile and be surrounded by golden rays of sunlight. The image should also include other symbols from the game, such as the Book of Helios symbol and the expandable symbol. Make sure the overall color scheme is bright and eye-catching, and that the image is of high-quality. Thanks!'\n\n$reviewTextRange = $reviewPara.Range.Duplicate\n$reviewTextRange.Find.Execute($oldReviewText) | Out-Null\n# Assigning .Text on the matched range replaces just that text while\n# keeping the paragraph's existing run/formatting structure (italic) intact.\n$reviewTextRange.Text = $dallePrompt\n"}
